$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B11 text from "Proveedor de lentes" to "Proveedor"
$ws.Range("B11").Value = "Proveedor"

# Delete row 12 ("7, Proveedor de monturas"), shifting the trailing
# empty-style row up to become the new row 12
$ws.Rows("12").Delete()

# Update the selected cell to match the post-edit state
$ws.Range("G13").Select()
